$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1) - column F "想去人数" updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 63
$wsExpo.Range("F8").Value = 95
$wsExpo.Range("F9").Value = 8624
$wsExpo.Range("F10").Value = 801
$wsExpo.Range("F11").Value = 325
$wsExpo.Range("F12").Value = 1140
$wsExpo.Range("F13").Value = 951
$wsExpo.Range("F14").Value = 98
$wsExpo.Range("F16").Value = 4
$wsExpo.Range("F17").Value = 230
$wsExpo.Range("F18").Value = 235
$wsExpo.Range("F21").Value = 993

# "全部类型" sheet (sheet4) - column F "想去人数" updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 63
$wsAll.Range("F10").Value = 95
$wsAll.Range("F11").Value = 8624
$wsAll.Range("F12").Value = 801
$wsAll.Range("F13").Value = 325
$wsAll.Range("F14").Value = 1140
$wsAll.Range("F15").Value = 951
$wsAll.Range("F16").Value = 98
$wsAll.Range("F18").Value = 4
$wsAll.Range("F19").Value = 230
$wsAll.Range("F20").Value = 235
$wsAll.Range("F21").Value = 61
$wsAll.Range("F23").Value = 993
